# Update crypto price/volume table with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.972.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.778.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.540.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.967.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0706"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.420.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  +3.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.529"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.692.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("D51").Style = "Normal"
